# BandColorConvention.xlsx - "Consolidated colors, reordered inputs on
# acoustic setup, Fixed positions"
#
# The sheet lists band members, their color, stage position and
# abbreviation. This commit renames/consolidates the stage-position
# vocabulary from a generic "stage left/right/center" scheme to a
# "downstage/upstage" scheme, and widens/adds a couple of columns so the
# new (longer) labels are readable.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-label the "Position" / "Abbreviated" columns (C, D) -------------
# Row 2 (Chad)    : Upstage Center / Drum  -> Upstage Center / USC
# Row 3 (John)    : Stage Right / SR              -> Downstage Right / DSR
# Row 4 (Jessica) : Center Stage Right / CSR      -> Downstage Center Right / DCR
# Row 5 (Paul)    : Center Stage Left / CSL       -> Downstage Center Left / DCL
# Row 6 (Kenzi)   : Stage Left / SL               -> Downstage Left / DSL
$ws.Range("D2").Value = "USC"

$ws.Range("C3").Value = "Downstage Right"
$ws.Range("D3").Value = "DSR"

$ws.Range("C4").Value = "Downstage Center Right"
$ws.Range("D4").Value = "DCR"

$ws.Range("C5").Value = "Downstage Center Left"
$ws.Range("D5").Value = "DCL"

$ws.Range("C6").Value = "Downstage Left"
$ws.Range("D6").Value = "DSL"

# --- Column widths --------------------------------------------------
# New columns A and B get an explicit best-fit-style width, and column C
# is widened to fit the longer "Downstage Center Right/Left" labels.
# (Column widths set through Excel's object model are stored in whole
# "standard character" units, so these are the closest values the
# COM layer can produce to the fitted widths.)
$ws.Columns.Item(1).ColumnWidth = 7.833333333333333
$ws.Columns.Item(2).ColumnWidth = 7.666666666666667
$ws.Columns.Item(3).ColumnWidth = 21.833333333333336

# --- Selection moved as part of the edit -----------------------------
$ws.Range("F9").Select() | Out-Null
